$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.816.08'
$ws.Cells.Item(2, 5).Value = '  +0.02%  '
$ws.Cells.Item(3, 4).Value = '1.633.41'
$ws.Cells.Item(3, 5).Value = '  +0.22%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 5).Value = '  -0.38%  '
$ws.Cells.Item(6, 5).Value = '  -0.47%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 5).Value = '  -0.23%  '
$ws.Cells.Item(9, 5).Value = '  -0.31%  '
$ws.Cells.Item(10, 4).Value = "'19.87"
$ws.Cells.Item(10, 5).Value = '  +2.22%  '
$ws.Cells.Item(11, 5).Value = '  +0.04%  '
$ws.Cells.Item(12, 4).Value = '1.661.29'
$ws.Cells.Item(12, 5).Value = '  +2.05%  '
$ws.Cells.Item(13, 5).Value = '  -0.35%  '
$ws.Cells.Item(14, 5).Value = '  +0.21%  '
$ws.Cells.Item(15, 5).Value = '  -0.24%  '
$ws.Cells.Item(16, 5).Value = '  +1.49%  '
$ws.Cells.Item(17, 4).Value = "'63.02"
$ws.Cells.Item(17, 5).Value = '  -0.29%  '
$ws.Cells.Item(18, 4).Value = '25.823.60'
$ws.Cells.Item(18, 5).Value = '  -0.04%  '
$ws.Cells.Item(20, 4).Value = "'193.74"
$ws.Cells.Item(20, 5).Value = '  -0.26%  '
$ws.Cells.Item(21, 5).Value = '  +1.24%  '
$ws.Cells.Item(22, 5).Value = '  +0.94%  '
$ws.Cells.Item(23, 4).Value = "'6.17"
$ws.Cells.Item(23, 5).Value = '  +2.53%  '
$ws.Cells.Item(24, 5).Value = '  -0.01%  '
$ws.Cells.Item(25, 5).Value = '  -2.51%  '
$ws.Cells.Item(26, 4).Value = "'139.42"
$ws.Cells.Item(26, 5).Value = '  -1.32%  '
$ws.Cells.Item(27, 5).Value = '  -2.99%  '
$ws.Cells.Item(28, 5).Value = '  +1.34%  '
$ws.Cells.Item(29, 5).Value = '  +0.60%  '
$ws.Cells.Item(30, 5).Value = '  +0.11%  '
$ws.Cells.Item(31, 4).Value = "'0.0494"
$ws.Cells.Item(31, 5).Value = '  +1.29%  '
$ws.Cells.Item(32, 5).Value = '  +0.73%  '
$ws.Cells.Item(33, 4).Value = "'3.24"
$ws.Cells.Item(33, 5).Value = '  +1.45%  '
$ws.Cells.Item(34, 5).Value = '  +1.05%  '
$ws.Cells.Item(35, 4).Value = "'2.39"
$ws.Cells.Item(35, 5).Value = '  +0.44%  '
$ws.Cells.Item(36, 4).Value = "'0.901"
$ws.Cells.Item(36, 5).Value = '  +0.58%  '
$ws.Cells.Item(37, 5).Value = '  +0.20%  '
$ws.Cells.Item(38, 4).Value = "'0.550"
$ws.Cells.Item(38, 5).Value = '  +0.55%  '
$ws.Cells.Item(39, 4).Value = '1.120.29'
$ws.Cells.Item(39, 5).Value = '  -0.97%  '
$ws.Cells.Item(40, 5).Value = '  +0.22%  '
$ws.Cells.Item(41, 5).Value = '  +0.76%  '
$ws.Cells.Item(42, 5).Value = '  -0.71%  '
$ws.Cells.Item(43, 4).Value = "'99.52"
$ws.Cells.Item(43, 5).Value = '  +2.10%  '
$ws.Cells.Item(44, 4).Value = "'0.800"
$ws.Cells.Item(44, 5).Value = '  +0.31%  '
$ws.Cells.Item(45, 4).Value = '0.0₆0107'
$ws.Cells.Item(45, 5).Value = '  -3.79%  '
$ws.Cells.Item(46, 4).Value = "'55.43"
$ws.Cells.Item(46, 5).Value = '  +0.33%  '
$ws.Cells.Item(47, 4).Value = "'0.420"
$ws.Cells.Item(47, 5).Value = '  -5.21%  '
$ws.Cells.Item(48, 5).Value = '  +0.30%  '
$ws.Cells.Item(49, 5).Value = '  -0.57%  '
$ws.Cells.Item(50, 4).Value = "'2.34"
$ws.Cells.Item(50, 5).Value = '  +6.02%  '
$ws.Cells.Item(51, 5).Value = '  -0.42%  '
